# worked on exporting time series files.
# Update the "Actual" hours tracking column (I) on the planning sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I header was "Actual (days)" -> now tracked in hours.
$ws.Range("I1").Value = "Actual (hours)"

# Actual time spent, revised from days to hours for several rows.
$ws.Range("I3").Value = 30
$ws.Range("I4").Value = 12
$ws.Range("I21").Value = 4

# New actual-hours entry for row 30 (previously blank).
$ws.Range("I30").Value = 1

# Leave the cell selection where the user was last working.
$ws.Range("I3").Select()
